$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block: card holder name
$ws.Range("C2").Value = "Hartmut"

# Card number (B3) must stay a TEXT value even though it's all digits.
# A plain Range.Value assignment of an all-digit string is auto-coerced
# to a number by this engine (and a 16-digit integer then loses precision
# under the General number format), so build it as a quoted-string formula
# in a scratch cell (a string-literal formula result is never reinterpreted
# as a number) and copy/paste-values it into B3. That preserves B3's exact
# original style (s="8") and type (text) with no other side effects, and
# unlike NumberFormat="@"/a leading apostrophe it doesn't leave behind an
# unused new style record.
$scratch = $ws.Range("G1")
$scratch.Formula = '="2570314725427075"'
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163, 0, $false, $false)  # xlPasteValues
$scratch.ClearContents()

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 04.04.2025"

# Transaction row 6
$ws.Range("B6").Value = "05.04."
$ws.Range("C6").Value = "06.04."
$ws.Range("D6").Value = "KARTENZ./05.04 EDEKA RO"
$ws.Range("E6").Value = "55,03-"

# Transaction row 7
$ws.Range("B7").Value = "07.04."
$ws.Range("C7").Value = "08.04."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,26-"

# Transaction row 8
$ws.Range("B8").Value = "10.04."
$ws.Range("C8").Value = "11.04."
$ws.Range("D8").Value = "KARTENZ./10.04 REWE RO"
$ws.Range("E8").Value = "92,45-"

# Transaction row 9
$ws.Range("B9").Value = "14.04."
$ws.Range("C9").Value = "15.04."
$ws.Range("D9").Value = "BURGER KING Dachau"
$ws.Range("E9").Value = "18,87-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 19.04.2025"
$ws.Range("E12").Value = "191,61-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 24.04.2025"
